$d = $word.ActiveDocument

# Replace "Table 7" with "Table 2"
$d.Content.Find.Execute("Table 7", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Table 2", 2)

# Replace every occurrence of "PTA" with "TA" (covers "PTA", "PTAs", etc.)
$range = $d.Content
$range.Find.Execute("PTA", $true, $false, $false, $false, $false,
                     $true, 1, $false, "TA", 2)
while ($range.Find.Found) {
    $range.Collapse(0)
    $range.End = $d.Content.End
    $range.Find.Execute("PTA", $true, $false, $false, $false, $false,
                         $true, 1, $false, "TA", 2)
}
